$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 43 with a new journal entry (13 April 2018, new Maven entry, 1 hour)
$newDate = (Get-Date -Year 2018 -Month 4 -Day 13 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A43").Value = $newDate
$ws.Range("B43").Value = "Résolution des problèmes liés à une mauvaise utilisation de Maven. Problèmes de dépendances, de version de java, etc."
$ws.Range("C43").Value = 1

# The row now holds wrapped text like the neighbouring entries, so give it the
# same taller row height (30pt) used by similar rows.
$ws.Rows.Item(43).RowHeight = 30

# Update the current selection to C44, matching where the user ended up after
# typing the new entry.
$ws.Range("C44").Select()
